{"js": "// Remove the two trailing \"page footer\" paragraphs that were dropped from\n// the rebuilt site (\"Ver no Jupiter Salvar em pdf Salvar em docx\" and the\n// \"\u00a9 2020 . Contact: luizeleno@usp.br. ...\" copyright line), leaving the\n// surrounding blank paragraphs untouched.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst targets = [\n  \"Ver no Jupiter Salvar em pdf Salvar em docx\",\n  \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n];\n\nfor (const paragraph of paragraphs.items) {\n  const text = paragraph.text.trim();\n  if (targets.indexOf(text) !== -1) {\n    paragraph.delete();\n  }\n}\n\nawait context.sync();\n", "ps1": "# Remove the two trailing \"page footer\" paragraphs that were dropped from\n# the rebuilt site (\"Ver no Jupiter Salvar em pdf Salvar em docx\" and the\n# \"\u00a9 2020 . Contact: luizeleno@usp.br. ...\" copyright line), leaving the\n# surrounding blank paragraphs untouched.\n$d = $word.ActiveDocument\n\n$needles = @(\n    \"Ver no Jupiter Salvar em pdf Salvar em docx\",\n    \"Contact: luizeleno@usp.br\"\n)\n\nforeach ($needle in $needles) {\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $found = $rng.Find.Execute($needle)\n    if ($found) {\n        $rng.Expand(4) | Out-Null   # wdParagraph -> whole paragraph incl. mark\n        $rng.Delete()\n    }\n}\n"}
